# update scripts wuth new tpm
# Recomputed NATMI ligand-receptor (Bmp7-Acvr2a) edge statistics.
# - Existing rows 2-9 (senders FAPs/ECs x targets FAPs/MuSCs/ECs/Resolving-Mac)
#   get refreshed TPM-derived values.
# - A third sender cluster, MuSCs, is added (rows 10-13), extending the table
#   to 3 senders x 4 targets = 12 data rows (header + rows 2-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp7"
$ws.Range("C2").Value = "Acvr2a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03443933333333333
$ws.Range("H2").Value = 0.103318
$ws.Range("I2").Value = 0.05823261822459219
$ws.Range("J2").Value = 0.0582326182245922
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 0.560963267195111
$ws.Range("R2").Value = 5.048669404756
$ws.Range("S2").Value = 0.0126766863031922
$ws.Range("T2").Value = 0.0126766863031922

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp7"
$ws.Range("C3").Value = "Acvr2a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03443933333333333
$ws.Range("H3").Value = 0.103318
$ws.Range("I3").Value = 0.05823261822459219
$ws.Range("J3").Value = 0.0582326182245922
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83271999999999
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 0.9509012183288887
$ws.Range("R3").Value = 8.558110964959999
$ws.Range("S3").Value = 0.02148853080942633
$ws.Range("T3").Value = 0.02148853080942634

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Acvr2a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03443933333333333
$ws.Range("H4").Value = 0.103318
$ws.Range("I4").Value = 0.05823261822459219
$ws.Range("J4").Value = 0.0582326182245922
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 0.9045947221166665
$ws.Range("R4").Value = 8.141352499049999
$ws.Range("S4").Value = 0.02044209343890573
$ws.Range("T4").Value = 0.02044209343890573

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("C5").Value = "Acvr2a"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03443933333333333
$ws.Range("H5").Value = 0.103318
$ws.Range("I5").Value = 0.05823261822459219
$ws.Range("J5").Value = 0.0582326182245922
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 0.1604255550884444
$ws.Range("R5").Value = 1.443829995796
$ws.Range("S5").Value = 0.003625307673067926
$ws.Range("T5").Value = 0.003625307673067927

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp7"
$ws.Range("C6").Value = "Acvr2a"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.51625
$ws.Range("H6").Value = 1.54875
$ws.Range("I6").Value = 0.8729143757654733
$ws.Range("J6").Value = 0.8729143757654734
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.28844733333333
$ws.Range("N6").Value = 48.865342
$ws.Range("O6").Value = 0.2176904746803693
$ws.Range("P6").Value = 0.2176904746803693
$ws.Range("Q6").Value = 8.408910935833333
$ws.Range("R6").Value = 75.6801984225
$ws.Range("S6").Value = 0.1900251448157042
$ws.Range("T6").Value = 0.1900251448157042

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp7"
$ws.Range("C7").Value = "Acvr2a"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.51625
$ws.Range("H7").Value = 1.54875
$ws.Range("I7").Value = 0.8729143757654733
$ws.Range("J7").Value = 0.8729143757654734
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.61090666666666
$ws.Range("N7").Value = 82.83271999999999
$ws.Range("O7").Value = 0.3690119294748028
$ws.Range("P7").Value = 0.3690119294748029
$ws.Range("Q7").Value = 14.25413056666667
$ws.Range("R7").Value = 128.2871751
$ws.Range("S7").Value = 0.3221158180675103
$ws.Range("T7").Value = 0.3221158180675104

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Bmp7"
$ws.Range("C8").Value = "Acvr2a"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.51625
$ws.Range("H8").Value = 1.54875
$ws.Range("I8").Value = 0.8729143757654733
$ws.Range("J8").Value = 0.8729143757654734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.266325
$ws.Range("N8").Value = 78.798975
$ws.Range("O8").Value = 0.3510419771967738
$ws.Range("P8").Value = 0.3510419771967739
$ws.Range("Q8").Value = 13.55999028125
$ws.Range("R8").Value = 122.03991253125
$ws.Range("S8").Value = 0.3064295883921993
$ws.Range("T8").Value = 0.3064295883921994

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Bmp7"
$ws.Range("C9").Value = "Acvr2a"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.51625
$ws.Range("H9").Value = 1.54875
$ws.Range("I9").Value = 0.8729143757654733
$ws.Range("J9").Value = 0.8729143757654734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.658207333333333
$ws.Range("N9").Value = 13.974622
$ws.Range("O9").Value = 0.06225561864805391
$ws.Range("P9").Value = 0.06225561864805392
$ws.Range("Q9").Value = 2.404799535833333
$ws.Range("R9").Value = 21.6431958225
$ws.Range("S9").Value = 0.05434382449005934
$ws.Range("T9").Value = 0.05434382449005935

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Bmp7"
$ws.Range("C10").Value = "Acvr2a"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.04072033333333334
$ws.Range("H10").Value = 0.122161
$ws.Range("I10").Value = 0.06885300600993445
$ws.Range("J10").Value = 0.06885300600993445
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 16.28844733333333
$ws.Range("N10").Value = 48.865342
$ws.Range("O10").Value = 0.2176904746803693
$ws.Range("P10").Value = 0.2176904746803693
$ws.Range("Q10").Value = 0.6632710048957778
$ws.Range("R10").Value = 5.969439044062
$ws.Range("S10").Value = 0.01498864356147295
$ws.Range("T10").Value = 0.01498864356147295

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Bmp7"
$ws.Range("C11").Value = "Acvr2a"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.04072033333333334
$ws.Range("H11").Value = 0.122161
$ws.Range("I11").Value = 0.06885300600993445
$ws.Range("J11").Value = 0.06885300600993445
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 27.61090666666666
$ws.Range("N11").Value = 82.83271999999999
$ws.Range("O11").Value = 0.3690119294748028
$ws.Range("P11").Value = 0.3690119294748029
$ws.Range("Q11").Value = 1.124325323102222
$ws.Range("R11").Value = 10.11892790792
$ws.Range("S11").Value = 0.02540758059786611
$ws.Range("T11").Value = 0.02540758059786611

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Bmp7"
$ws.Range("C12").Value = "Acvr2a"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04072033333333334
$ws.Range("H12").Value = 0.122161
$ws.Range("I12").Value = 0.06885300600993445
$ws.Range("J12").Value = 0.06885300600993445
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 26.266325
$ws.Range("N12").Value = 78.798975
$ws.Range("O12").Value = 0.3510419771967738
$ws.Range("P12").Value = 0.3510419771967739
$ws.Range("Q12").Value = 1.069573509441667
$ws.Range("R12").Value = 9.626161584975
$ws.Range("S12").Value = 0.02417029536566874
$ws.Range("T12").Value = 0.02417029536566874

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Bmp7"
$ws.Range("C13").Value = "Acvr2a"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04072033333333334
$ws.Range("H13").Value = 0.122161
$ws.Range("I13").Value = 0.06885300600993445
$ws.Range("J13").Value = 0.06885300600993445
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.658207333333333
$ws.Range("N13").Value = 13.974622
$ws.Range("O13").Value = 0.06225561864805391
$ws.Range("P13").Value = 0.06225561864805392
$ws.Range("Q13").Value = 0.1896837553491111
$ws.Range("R13").Value = 1.707153798142
$ws.Range("S13").Value = 0.004286486484926643
$ws.Range("T13").Value = 0.004286486484926644
